$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 26800
$ws.Range("C3").Value = 52500
$ws.Range("D3").Value = 71700
$ws.Range("E3").Value = 85900
$ws.Range("B4").Value = 110
$ws.Range("C4").Value = 215
$ws.Range("D4").Value = 294
$ws.Range("E4").Value = 352
$ws.Range("B5").Value = 37.03445
$ws.Range("C5").Value = 37.7654
$ws.Range("D5").Value = 39.85357
$ws.Range("E5").Value = 41.63467
$ws.Range("B11").Value = 6965
$ws.Range("C11").Value = 11400
$ws.Range("D11").Value = 18500
$ws.Range("E11").Value = 20500
$ws.Range("B12").Value = 913
$ws.Range("C12").Value = 1491
$ws.Range("D12").Value = 2426
$ws.Range("E12").Value = 2684
$ws.Range("B13").Value = 138.4
$ws.Range("C13").Value = 158.48
$ws.Range("D13").Value = 155.77
$ws.Range("E13").Value = 173.23
$ws.Range("B19").Value = 49600
$ws.Range("C19").Value = 70800
$ws.Range("D19").Value = 83200
$ws.Range("E19").Value = 101000
$ws.Range("B20").Value = 203
$ws.Range("C20").Value = 290
$ws.Range("D20").Value = 341
$ws.Range("E20").Value = 415
$ws.Range("B21").Value = 18.1
$ws.Range("C21").Value = 25.4
$ws.Range("D21").Value = 30.87
$ws.Range("E21").Value = 33.34
$ws.Range("B27").Value = 3878
$ws.Range("C27").Value = 6005
$ws.Range("D27").Value = 8370
$ws.Range("E27").Value = 9122
$ws.Range("B28").Value = 508
$ws.Range("C28").Value = 787
$ws.Range("D28").Value = 1097
$ws.Range("E28").Value = 1196
$ws.Range("B29").Value = 201.13
$ws.Range("C29").Value = 249.36
$ws.Range("D29").Value = 243.57
$ws.Range("E29").Value = 275.26
$ws.Range("B35").Value = 218000
$ws.Range("C35").Value = 379000
$ws.Range("D35").Value = 512000
$ws.Range("E35").Value = 615000
$ws.Range("B36").Value = 895
$ws.Range("C36").Value = 1552
$ws.Range("D36").Value = 2097
$ws.Range("E36").Value = 2521
$ws.Range("B37").Value = 4.234970000000001
$ws.Range("C37").Value = 4.90655
$ws.Range("D37").Value = 5.53202
$ws.Range("E37").Value = 6.1661
$ws.Range("B43").Value = 6649
$ws.Range("C43").Value = 11900
$ws.Range("D43").Value = 16300
$ws.Range("E43").Value = 19300
$ws.Range("B44").Value = 872
$ws.Range("C44").Value = 1561
$ws.Range("D44").Value = 2142
$ws.Range("E44").Value = 2532
$ws.Range("B45").Value = 145.37
$ws.Range("C45").Value = 158.42
$ws.Range("D45").Value = 177.88
$ws.Range("E45").Value = 198.09
$ws.Range("B51").Value = 49900
$ws.Range("C51").Value = 61200
$ws.Range("D51").Value = 97500
$ws.Range("E51").Value = 113000
$ws.Range("B52").Value = 204
$ws.Range("C52").Value = 251
$ws.Range("D52").Value = 399
$ws.Range("E52").Value = 465
$ws.Range("B53").Value = 17.48
$ws.Range("C53").Value = 25.76
$ws.Range("D53").Value = 25.78
$ws.Range("E53").Value = 28.35
$ws.Range("B59").Value = 3864
$ws.Range("C59").Value = 5769
$ws.Range("D59").Value = 7094
$ws.Range("E59").Value = 7699
$ws.Range("B60").Value = 506
$ws.Range("C60").Value = 756
$ws.Range("D60").Value = 930
$ws.Range("E60").Value = 1009
$ws.Range("B61").Value = 186.32
$ws.Range("C61").Value = 227.78
$ws.Range("D61").Value = 262.4
$ws.Range("E61").Value = 285
$ws.Range("B67").Value = 33700
$ws.Range("C67").Value = 65600
$ws.Range("D67").Value = 97100
$ws.Range("E67").Value = 120000
$ws.Range("B68").Value = 138
$ws.Range("C68").Value = 269
$ws.Range("D68").Value = 398
$ws.Range("E68").Value = 493
$ws.Range("B69").Value = 29.29052
$ws.Range("C69").Value = 29.95759
$ws.Range("D69").Value = 30.40274
$ws.Range("E69").Value = 31.45437
$ws.Range("B75").Value = 8982
$ws.Range("C75").Value = 17400
$ws.Range("D75").Value = 23600
$ws.Range("E75").Value = 28800
$ws.Range("B76").Value = 1177
$ws.Range("C76").Value = 2275
$ws.Range("D76").Value = 3097
$ws.Range("E76").Value = 3781
$ws.Range("B77").Value = 102.78
$ws.Range("C77").Value = 106.32
$ws.Range("D77").Value = 120.81
$ws.Range("E77").Value = 129.47
$ws.Range("B83").Value = 79500
$ws.Range("C83").Value = 122000
$ws.Range("D83").Value = 151000
$ws.Range("E83").Value = 165000
$ws.Range("B84").Value = 326
$ws.Range("C84").Value = 500
$ws.Range("D84").Value = 618
$ws.Range("E84").Value = 675
$ws.Range("B85").Value = 10.71
$ws.Range("C85").Value = 13.7
$ws.Range("D85").Value = 15.94
$ws.Range("E85").Value = 19.13
$ws.Range("B91").Value = 6059
$ws.Range("C91").Value = 8677
$ws.Range("D91").Value = 10900
$ws.Range("E91").Value = 11600
$ws.Range("B92").Value = 794
$ws.Range("C92").Value = 1137
$ws.Range("D92").Value = 1423
$ws.Range("E92").Value = 1517
$ws.Range("B93").Value = 117.46
$ws.Range("C93").Value = 141.39
$ws.Range("D93").Value = 161.64
$ws.Range("E93").Value = 211.02
$ws.Range("B99").Value = 315000
$ws.Range("C99").Value = 643000
$ws.Range("D99").Value = 793000
$ws.Range("E99").Value = 814000
$ws.Range("B100").Value = 1291
$ws.Range("C100").Value = 2632
$ws.Range("D100").Value = 3247
$ws.Range("E100").Value = 3335
$ws.Range("B101").Value = 2.92332
$ws.Range("C101").Value = 2.85878
$ws.Range("D101").Value = 3.40918
$ws.Range("E101").Value = 4.44653
$ws.Range("B107").Value = 10600
$ws.Range("C107").Value = 18600
$ws.Range("D107").Value = 25400
$ws.Range("E107").Value = 27900
$ws.Range("B108").Value = 1384
$ws.Range("C108").Value = 2440
$ws.Range("D108").Value = 3328
$ws.Range("E108").Value = 3652
$ws.Range("B109").Value = 90.12
$ws.Range("C109").Value = 101.93
$ws.Range("D109").Value = 111.67
$ws.Range("E109").Value = 130.68
$ws.Range("B115").Value = 83400
$ws.Range("C115").Value = 127000
$ws.Range("D115").Value = 141000
$ws.Range("E115").Value = 170000
$ws.Range("B116").Value = 342
$ws.Range("C116").Value = 518
$ws.Range("D116").Value = 578
$ws.Range("E116").Value = 697
$ws.Range("B117").Value = 9.619999999999999
$ws.Range("C117").Value = 11.98
$ws.Range("D117").Value = 16.32
$ws.Range("E117").Value = 17.57
$ws.Range("B123").Value = 5224
$ws.Range("C123").Value = 6918
$ws.Range("D123").Value = 8439
$ws.Range("E123").Value = 9041
$ws.Range("B124").Value = 685
$ws.Range("C124").Value = 907
$ws.Range("D124").Value = 1106
$ws.Range("E124").Value = 1185
$ws.Range("B125").Value = 116.64
$ws.Range("C125").Value = 148.98
$ws.Range("D125").Value = 178.42
$ws.Range("E125").Value = 229.1
